# Outstandings.xlsx edit
#
# The only substantive change is on the "Sale 22-23" sheet: the three
# Putzmeister rows 8:10 (b23-24AQ307 / b23-24MQ311 / b23-24AQ312 invoices)
# are removed as whole rows, which shifts every row below them up by three
# and lets Excel auto-adjust the relative formulas that referenced the
# shifted cells. A running-total formula is then written into the now-empty
# F7 cell. Finally the active sheet/selection bookkeeping is updated so that
# "Purchase 22-23" (sheet 1) is the selected tab with F17 selected, and the
# "Sale 22-23" sheet keeps a B5:F7 selection instead of being the active tab.

$wb = $excel.ActiveWorkbook

$wsSale = $wb.Worksheets.Item("Sale 22-23")
$wsPurchase = $wb.Worksheets.Item("Purchase 22-23")

# Remove the three rows for the deleted invoices (b23-24AQ307, b23-24MQ311,
# b23-24AQ312); everything below shifts up and formulas re-point themselves.
$wsSale.Rows("8:10").Delete()

# Add the new subtotal formula that now belongs in the (still blank) F7 cell.
$wsSale.Range("F7").Formula = "=E5+E6+E7"

# Update the view/selection state: Sale 22-23 keeps a non-active selection ...
$wsSale.Range("B5:F7").Select()

# ... while Purchase 22-23 becomes the active/selected sheet with F17 picked.
$wsPurchase.Activate()
$wsPurchase.Range("F17").Select()
